$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 55 (shifts rows 55..94 down to 56..95,
# and copies formatting from the row above automatically).
$ws.Rows("55:55").Insert()

# Fill in the new product row (OXALEPTAL 300 MG 30 F.C.TABS.)
$ws.Cells.Item(55, 1).Value = 49
$ws.Cells.Item(55, 3).Value = "OXALEPTAL 300 MG 30 F.C.TABS."
$ws.Cells.Item(55, 8).Value = "1:0"
$ws.Cells.Item(55, 12).Value = 1
$ws.Cells.Item(55, 14).Value = "180.00"
$ws.Cells.Item(55, 16).Value = 59.4
$ws.Cells.Item(55, 17).Value = "0:1"

# Renumber the sequence column (A) for every following data row (56..93)
for ($r = 56; $r -le 93; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 7
}

# Update the grand total (now on row 94, column P)
$ws.Cells.Item(94, 16).Value = 5396.685

# Update the generated timestamp string
$ws.Cells.Item(95, 1).Value = "Saturday, 27 September, 2025 8:52 PM"
